$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 94 (pushes existing rows 94.. down by one,
# carrying the date-number style already present on column D down with them).
$ws.Rows.Item(94).Insert()

# Populate the newly inserted row 94 with a new weekly price record
# (same attributes as the following record, new date 2023-08-08 = serial 45146).
$ws.Cells.Item(94, 1).Value = 5
$ws.Cells.Item(94, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(94, 3).Value = "Maule"
$ws.Cells.Item(94, 4).Value = 45146
$ws.Cells.Item(94, 5).Value = 7
$ws.Cells.Item(94, 6).Value = 100112013
$ws.Cells.Item(94, 7).Value = "Alcachofa"
$ws.Cells.Item(94, 8).Value = "Madrigal"
$ws.Cells.Item(94, 9).Value = "Primera"
$ws.Cells.Item(94, 10).Value = 300
$ws.Cells.Item(94, 11).Value = 13000
$ws.Cells.Item(94, 12).Value = 13000
$ws.Cells.Item(94, 13).Value = 13000
$ws.Cells.Item(94, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(94, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(94, 16).Value = 325
$ws.Cells.Item(94, 17).Value = 40
$ws.Cells.Item(94, 18).Value = "Hortaliza"
